# converter jin (斤) to kg unit
# Updates the "重量(公斤)" (weight in kg) column (N) on the active sheet,
# converting values that were recorded in jin/other units (e.g. "25.8斤",
# "90隻") into their kilogram-equivalent numeric text, and clears the two
# "(庫存)" placeholder cells that no longer apply.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells must keep a *text* data type (the source values were never real
# numbers, e.g. "7.199999999999999"/"15.0" as literal strings) so force
# the number format to Text first, otherwise Excel auto-coerces a
# numeric-looking string into a true number. ClearFormats() afterwards
# drops the temporary Text number-format again so no stray style survives
# on the cells once the values are committed.
$weightRange = $ws.Range("N2:N10")
$weightRange.NumberFormat = "@"

$ws.Range("N2").Value  = "15.0"                 # 白米: 25.8斤 -> 15.0 kg
$ws.Range("N3").Value  = "7.199999999999999"    # 檸檬雞翅: 90隻 -> 7.199999999999999 kg
$ws.Range("N4").Value  = "2.4"                  # 金針菇: 4斤 -> 2.4 kg
$ws.Range("N5").Value  = "1.7999999999999998"   # 杏鮑菇: 3斤 -> 1.7999999999999998 kg
$ws.Range("N6").Value  = "1.7999999999999998"   # 袖珍菇: 3斤 -> 1.7999999999999998 kg
$ws.Range("N7").Value  = "4.8"                  # 小白菜: 8斤 -> 4.8 kg
$ws.Range("N8").Value  = ""                     # 福菜: (庫存) -> (cleared)
$ws.Range("N9").Value  = ""                     # 竹筍: (庫存) -> (cleared)
$ws.Range("N10").Value = "2.4"                  # 排骨: 4斤 -> 2.4 kg

$weightRange.ClearFormats()
